$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Board")
$ws.Range("B2").Value = "Set up Python backend with Flask"
$ws.Range("B3").Value = "Set up React frontend with Vite, Tailwind CSS, and responsive navigation"
$ws.Range("B4").Value = "Connect backend to SQLite and verify local connection"
$ws.Range("B8").Value = "Create .env configuration for sensitive keys"
